$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheet name: $($ws.Name)"
Write-Host ("A54: " + $ws.Cells.Item(54,1).Value)
